$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price-like text values in column D keep their exact textual
# representation (leading/trailing zeros) instead of being auto-converted
# to numbers by Excel when assigned through .Value.
$priceCells = @{
    "D2" = "246.10"
    "D3" = "22.17"
    "D4" = "5.369"
    "D5" = "0.05867"
    "D6" = "3.384"
    "D8" = "0.8130"
    "D9" = "0.9651"
    "D10" = "0.01118"
    "D11" = "0.1421"
    "D12" = "0.03640"
    "D13" = "0.07362"
    "D14" = "0.03016"
    "D15" = "4.475"
    "D16" = "0.09398"
    "D17" = "0.001590"
    "D18" = "0.04828"
    "D19" = "0.006246"
    "D20" = "0.004083"
    "D21" = "0.0009832"
    "D22" = "0.00009706"
    "D25" = "0.3252"
    "D27" = "0.0002473"
    "D40" = "0.03841"
    "D45" = "0.00005665"
    "D46" = "0.00000000751"
    "D47" = "0.6515"
    "D48" = "0.07749"
}
foreach ($addr in $priceCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$addr]
}

# Plain text cells (coin names, links, volume labels) -- these are not
# numeric-looking, so a direct .Value assignment keeps them as text.
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "MandalaExchangeToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("E13").Value = "12MandalaExchangeTokenMDX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitMartToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("E16").Value = "15BitMartTokenBMX"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
